$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update label text from "昨日" (yesterday) to "前日" (day before yesterday)
$ws.Range("G2").Value = "前日"
$ws.Range("H2").Value = "前日"

# Update selection to H2
$ws.Range("H2").Select()
